$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This change represents a localization "handoff" run: the old source file
# c44fd665-2a49-403a-99df-b60d2f6c7a10.md was superseded by a new source file
# 6ef09566-1f3e-4049-a389-5ba7897ab1dd.md (status flips from the old failure
# state "Handoff transform failed" to "Ready for handoff"), a brand-new
# source file ffff868d2531-c139-461c-a010-804db51013f8.md is now also being
# tracked (also "Ready for handoff"), and the unrelated .localization-config
# row simply shifts down to make room. The per-language sheets additionally
# gain handoff-file / handoff-datetime / handoff-reason detail for the two
# "Ready for handoff" rows.
# ---------------------------------------------------------------------------

$oldMd        = "c44fd665-2a49-403a-99df-b60d2f6c7a10.md"
$newMd        = "6ef09566-1f3e-4049-a389-5ba7897ab1dd.md"
$newMd2       = "ffff868d2531-c139-461c-a010-804db51013f8.md"
$cfgName      = ".localization-config"

$newMdUrl     = "https://github.com/OpenLocalizationTest/oltest/blob/92d123faba7748170c7859b78b8858d0bf204f00/e2e/$newMd"
$newMd2Url    = "https://github.com/OpenLocalizationTest/oltest/blob/92d123faba7748170c7859b78b8858d0bf204f00/e2e/$newMd2"
$cfgUrl       = "https://github.com/OpenLocalizationTest/oltest/blob/ae5fe9ee18c6dddb6ba3783d447e33120deccdcc/.localization-config"

$zhXlfName    = "6ef09566-1f3e-4049-a389-5ba7897ab1dd.92d123faba7748170c7859b78b8858d0bf204f00.zh-cn.xlf"
$deXlfName    = "6ef09566-1f3e-4049-a389-5ba7897ab1dd.92d123faba7748170c7859b78b8858d0bf204f00.de-de.xlf"
$zhXlfUrl     = "https://github.com/OpenLocalizationTest/oltest/blob/92d123faba7748170c7859b78b8858d0bf204f00/loc/$zhXlfName"
$deXlfUrl     = "https://github.com/OpenLocalizationTest/oltest/blob/92d123faba7748170c7859b78b8858d0bf204f00/loc/$deXlfName"

$zhHandoffDt  = "2016-01-19 07:15:24"
$deHandoffDt  = "2016-01-19 07:15:35"
$epoch        = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Hyperlinks.Delete()

$ov.Range("A2").Value = $newMd
$ov.Range("B2").Value = "Ready for handoff"
$ov.Range("C2").Value = "Ready for handoff"

$ov.Range("A3").Value = $newMd2
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"

$ov.Range("A4").Value = $cfgName
$ov.Range("B4").Value = "Not to be localized"
$ov.Range("C4").Value = "Not to be localized"

$ov.Hyperlinks.Add($ov.Range("A2"), $newMdUrl, "", "", $newMd)
$ov.Hyperlinks.Add($ov.Range("A3"), $newMd2Url, "", "", $newMd2)
$ov.Hyperlinks.Add($ov.Range("A4"), $cfgUrl, "", "", $cfgName)

# ---------------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | Status | Latest Handoff File |
#   Latest Handoff Datetime | Latest Target File | Latest Handback File |
#   Latest Handback DateTime | Handoff Reason | Dependency From
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Hyperlinks.Delete()

$zh.Range("A2").Value = $newMd
$zh.Range("B2").Value = "Ready for handoff"
$zh.Range("C2").Value = $zhXlfName
$zh.Range("D2").Value = $zhHandoffDt
$zh.Range("G2").Value = $epoch
$zh.Range("H2").Value = "Include"

$zh.Range("A3").Value = $newMd2
$zh.Range("B3").Value = "Ready for handoff"
$zh.Range("C3").Value = $zhXlfName
$zh.Range("D3").Value = $zhHandoffDt
$zh.Range("G3").Value = $epoch
$zh.Range("H3").Value = "Include"

$zh.Range("A4").Value = $cfgName
$zh.Range("B4").Value = "Not to be localized"
$zh.Range("D4").Value = $epoch
$zh.Range("G4").Value = $epoch
$zh.Range("H4").Value = "Ignored"

$zh.Hyperlinks.Add($zh.Range("A2"), $newMdUrl, "", "", $newMd)
$zh.Hyperlinks.Add($zh.Range("C2"), $zhXlfUrl, "", "", $zhXlfName)
$zh.Hyperlinks.Add($zh.Range("A3"), $newMd2Url, "", "", $newMd2)
$zh.Hyperlinks.Add($zh.Range("C3"), $zhXlfUrl, "", "", $zhXlfName)
$zh.Hyperlinks.Add($zh.Range("A4"), $cfgUrl, "", "", $cfgName)

# ---------------------------------------------------------------------------
# Sheet "de-de": same column layout as zh-cn, but German handoff file/time
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Hyperlinks.Delete()

$de.Range("A2").Value = $newMd
$de.Range("B2").Value = "Ready for handoff"
$de.Range("C2").Value = $deXlfName
$de.Range("D2").Value = $deHandoffDt
$de.Range("G2").Value = $epoch
$de.Range("H2").Value = "Include"

$de.Range("A3").Value = $newMd2
$de.Range("B3").Value = "Ready for handoff"
$de.Range("C3").Value = $deXlfName
$de.Range("D3").Value = $deHandoffDt
$de.Range("G3").Value = $epoch
$de.Range("H3").Value = "Include"

$de.Range("A4").Value = $cfgName
$de.Range("B4").Value = "Not to be localized"
$de.Range("D4").Value = $epoch
$de.Range("G4").Value = $epoch
$de.Range("H4").Value = "Ignored"

$de.Hyperlinks.Add($de.Range("A2"), $newMdUrl, "", "", $newMd)
$de.Hyperlinks.Add($de.Range("C2"), $deXlfUrl, "", "", $deXlfName)
$de.Hyperlinks.Add($de.Range("A3"), $newMd2Url, "", "", $newMd2)
$de.Hyperlinks.Add($de.Range("C3"), $deXlfUrl, "", "", $deXlfName)
$de.Hyperlinks.Add($de.Range("A4"), $cfgUrl, "", "", $cfgName)
